# Update the carjacking-by-neighborhood-by-month workbook with the
# 2021-09-19 data refresh (extends the "through" date from Sept 10 to
# Sept 11, and updates the partial-September counts for several
# neighborhoods).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet/tab name reflects the new "through" date.
$ws.Name = "Through 2021-09-11"

# Header cell (shared string) describing the partial current month.
$ws.Range("B1").Value = "September 2021 (through September 11)"

# Per-neighborhood September-2021-partial-month counts (column B) and a
# handful of other historical cells that were corrected/filled in.
$ws.Range("B2").Value = 8      # Garfield Park
$ws.Range("B3").Value = 3      # North Lawndale
$ws.Range("K3").Value = 6      # North Lawndale, September 2020
$ws.Range("T3").Value = 3      # North Lawndale, September 2019
$ws.Range("AC3").Value = 1     # North Lawndale, September 2018
$ws.Range("K4").Value = 4      # Humboldt Park, September 2020
$ws.Range("AC6").Value = 1     # Roseland, September 2018
$ws.Range("B7").Value = 4      # Auburn Gresham
$ws.Range("B9").Value = 2      # Little Village
$ws.Range("AC9").Value = 1     # Little Village, September 2018
$ws.Range("B10").Value = 1     # West Town
$ws.Range("B11").Value = 1     # Little Italy, UIC
$ws.Range("AC16").Value = 1    # West Pullman, September 2018
$ws.Range("AL17").Value = 3    # South Shore, September 2017
$ws.Range("K20").Value = 2     # Englewood, September 2020
$ws.Range("B27").Value = 1     # Avalon Park
$ws.Range("BD30").Value = 1    # Lincoln Park, July 2015
$ws.Range("AL32").Value = 1    # Chicago Lawn, September 2017
$ws.Range("AU38").Value = 1    # East Village, August 2018
$ws.Range("K49").Value = 1     # North Center, September 2020
$ws.Range("K55").Value = 2     # Grand Crossing, September 2020
$ws.Range("B97").Value = 1     # Uptown
$ws.Range("AC99").Value = 1    # West Ridge, September 2018
